$wb = $excel.ActiveWorkbook

# --- Rename the three PR_DOA cassette sheets, appending "_1" ---
$renames = @(
    @{ Old = "PR_DOA_5_Cassette"; New = "PR_DOA_5_Cassette_1" },
    @{ Old = "PR_DOA_4_Cassette"; New = "PR_DOA_4_Cassette_1" },
    @{ Old = "PR_DOA_3_Cassette"; New = "PR_DOA_3_Cassette_1" }
)

foreach ($r in $renames) {
    $ws = $wb.Worksheets.Item($r.Old)
    $ws.Name = $r.New

    # Keep the sheet-scoped Print_Area defined name's RefersTo text in sync
    # with the new sheet name (mirrors Excel's own rename behaviour).
    $dnName = "$($r.New)!Print_Area"
    $dn = $wb.Names.Item($dnName)
    $dn.RefersTo = "=$($r.New)!`$A`$1:`$H`$42"
}

# --- Switch the active tab to PR_DOA_3_Cassette_1 (was KET_Uncut_Sheet_1) ---
$target = $wb.Worksheets.Item("PR_DOA_3_Cassette_1")
$target.Activate()
$target.Range("C16").Select() | Out-Null
